$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 6386.1816
$ws.Range("I18").Value = 7037
$ws.Range("K18").Value = 7037
$ws.Range("M18").Value = -6753
$ws.Range("H40").Value = 1249.5
$ws.Range("I40").Value = 1149.5
$ws.Range("J40").Value = 1349.5
$ws.Range("K40").Value = 1149.5
$ws.Range("L40").Value = 1349.5
$ws.Range("M40").Value = -974.5
$ws.Range("N40").Value = -1699.5
$ws.Range("H103").Value = 487.25
$ws.Range("I103").Value = 499.66666
$ws.Range("J103").Value = 450
$ws.Range("K103").Value = 1498.99998
$ws.Range("L103").Value = 1350
$ws.Range("M103").Value = -912.9999800000001
$ws.Range("N103").Value = -2522
$ws.Range("H111").Value = 6691.75
$ws.Range("I111").Value = 392.69232
$ws.Range("K111").Value = 1178.07696
$ws.Range("M111").Value = 1888.92304

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2632.5278
$ws.Range("I61").Value = 1848.2593
$ws.Range("K61").Value = 1848.2593
$ws.Range("M61").Value = -1636.2593
$ws.Range("H110").Value = 5075.727
$ws.Range("I110").Value = 7553.067
$ws.Range("J110").Value = 3011.2778
$ws.Range("K110").Value = 7553.067
$ws.Range("L110").Value = 3011.2778
$ws.Range("M110").Value = -5508.067
$ws.Range("N110").Value = -7101.2778
$ws.Range("H136").Value = 2632.5278
$ws.Range("I136").Value = 1848.2593
$ws.Range("K136").Value = 5544.7779
$ws.Range("M136").Value = -2994.7779

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H101").Value = 45555
$ws.Range("J101").Value = 45555
$ws.Range("L101").Value = 45555
$ws.Range("N101").Value = -52045
$ws.Range("H105").Value = 5433.7334
$ws.Range("I105").Value = 5494.3335
$ws.Range("K105").Value = 5494.3335
$ws.Range("M105").Value = -3747.3335

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8291
$ws.Range("I16").Value = 6436.875
$ws.Range("J16").Value = 11999.25
$ws.Range("K16").Value = 6436.875
$ws.Range("L16").Value = 11999.25
$ws.Range("M16").Value = -6149.875
$ws.Range("N16").Value = -12573.25
$ws.Range("H113").Value = 8291
$ws.Range("I113").Value = 6436.875
$ws.Range("J113").Value = 11999.25
$ws.Range("K113").Value = 6436.875
$ws.Range("L113").Value = 11999.25
$ws.Range("M113").Value = -4266.875
$ws.Range("N113").Value = -16339.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1027.7368
$ws.Range("I5").Value = 904.1
$ws.Range("J5").Value = 1165.1111
$ws.Range("K5").Value = 2712.3
$ws.Range("L5").Value = 3495.3333
$ws.Range("M5").Value = -2600.3
$ws.Range("N5").Value = -3719.3333
$ws.Range("H21").Value = 321.66666
$ws.Range("I21").Value = 220
$ws.Range("J21").Value = 372.5
$ws.Range("K21").Value = 660
$ws.Range("L21").Value = 1117.5
$ws.Range("M21").Value = -487
$ws.Range("N21").Value = -1463.5
$ws.Range("H122").Value = 2201.5264
$ws.Range("J122").Value = 2485.4285
$ws.Range("L122").Value = 22368.8565
$ws.Range("N122").Value = -27268.8565
$ws.Range("H135").Value = 1027.7368
$ws.Range("I135").Value = 904.1
$ws.Range("J135").Value = 1165.1111
$ws.Range("K135").Value = 8136.900000000001
$ws.Range("L135").Value = 10485.9999
$ws.Range("M135").Value = -5601.900000000001
$ws.Range("N135").Value = -15555.9999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6547.7
$ws.Range("I80").Value = 3901.6667
$ws.Range("J80").Value = 7681.7144
$ws.Range("K80").Value = 3901.6667
$ws.Range("L80").Value = 7681.7144
$ws.Range("M80").Value = -2903.6667
$ws.Range("N80").Value = -9677.714400000001
$ws.Range("H83").Value = 6547.7
$ws.Range("I83").Value = 3901.6667
$ws.Range("J83").Value = 7681.7144
$ws.Range("K83").Value = 19508.3335
$ws.Range("L83").Value = 38408.572
$ws.Range("M83").Value = -14516.3335
$ws.Range("N83").Value = -48392.572
$ws.Range("H102").Value = 27778984
$ws.Range("I102").Value = 1133.3939
$ws.Range("K102").Value = 1133.3939
$ws.Range("M102").Value = 488.6061
$ws.Range("H107").Value = 1591.1875
$ws.Range("I107").Value = 1427.6923
$ws.Range("K107").Value = 1427.6923
$ws.Range("M107").Value = 492.3077000000001
$ws.Range("H113").Value = 1863.75
$ws.Range("I113").Value = 1863.75
$ws.Range("K113").Value = 1863.75
$ws.Range("M113").Value = 306.25
$ws.Range("H122").Value = 2269.875
$ws.Range("I122").Value = 2028.1666
$ws.Range("J122").Value = 2995
$ws.Range("K122").Value = 6084.4998
$ws.Range("L122").Value = 8985
$ws.Range("M122").Value = -3634.4998
$ws.Range("N122").Value = -13885
$ws.Range("H124").Value = 29999
$ws.Range("J124").Value = 29999
$ws.Range("L124").Value = 29999
$ws.Range("N124").Value = -39819
$ws.Range("H126").Value = 25266
$ws.Range("J126").Value = 4250
$ws.Range("L126").Value = 12750
$ws.Range("N126").Value = -17690

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7199.2
$ws.Range("I40").Value = 5698.6
$ws.Range("K40").Value = 5698.6
$ws.Range("M40").Value = -5562.6
$ws.Range("H122").Value = 7373.923
$ws.Range("I122").Value = 7126.8
$ws.Range("J122").Value = 8197.666999999999
$ws.Range("K122").Value = 21380.4
$ws.Range("L122").Value = 24593.001
$ws.Range("M122").Value = -18930.4
$ws.Range("N122").Value = -29493.001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 1000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459
$ws.Range("H113").Value = 1885.2
$ws.Range("I113").Value = 1710.3572
$ws.Range("J113").Value = 2293.1667
$ws.Range("K113").Value = 5131.071599999999
$ws.Range("L113").Value = 6879.500100000001
$ws.Range("M113").Value = -2961.071599999999
$ws.Range("N113").Value = -11219.5001
$ws.Range("H122").Value = 3005.3333
$ws.Range("I122").Value = 2763.25
$ws.Range("J122").Value = 3780
$ws.Range("K122").Value = 8289.75
$ws.Range("L122").Value = 11340
$ws.Range("M122").Value = -5839.75
$ws.Range("N122").Value = -16240
$ws.Range("H125").Value = 125026250
$ws.Range("J125").Value = 125026250
$ws.Range("L125").Value = 125026250
$ws.Range("N125").Value = -125036090
$ws.Range("H132").Value = 2127.1667
$ws.Range("I132").Value = 1979.2565
$ws.Range("K132").Value = 5937.7695
$ws.Range("M132").Value = -3407.7695
